$d = $word.ActiveDocument

$fields = @("to", "from", "date", "subject", "question", "answer", "facts", "discussion", "conclusion")

foreach ($f in $fields) {
    $old = "{{" + $f + "}}"
    $new = "{" + $f + "}"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
